# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gains a new (blank) column before the old
# column N ("Late"), pushing the old N/O/P ("Late"/"heading"/"Outstanding")
# columns one place to the right (O/P/Q), and the sheet becomes the active /
# selected sheet (instead of "Transactions").

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (14th column). This shifts the existing
# N, O, P columns (and their widths/styles) one column to the right,
# matching Excel's native "Insert Column" behaviour.
$wsSchedule.Columns.Item(14).Insert()

# Give the newly inserted column the same width as column M (10.7109375),
# mirroring the format that is carried over when a column is inserted next
# to an existing, explicitly-sized column.
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with P7 selected, and
# deactivate the previously active "Transactions" sheet.
$wsSchedule.Activate()
$wsSchedule.Range("P7").Select()
